# "keywords and locator helper"
#
# testcases sheet: row 5 ("employee creation") Run Mode flips from "yes" to "no"
# teststeps sheet: new C2 locator-helper cell (single space) on the
#                  "open desired browser" step row
# Also updates the remembered cell selection on the testcases sheet to B5
# (the cell that was just edited) while leaving teststeps as the active tab,
# matching the workbook's original view state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("testcases")
$ws2 = $wb.Worksheets.Item("teststeps")

# Run Mode for "employee creation" (row 5): yes -> no
$ws1.Range("B5").Value = "no"

# New Loc Type placeholder (single space) for the "open desired browser" step
$ws2.Range("C2").Value = " "

# Refresh remembered selections: testcases now points at B5, but keep
# teststeps as the visible/active sheet (as it was before the edit).
$ws1.Range("B5").Select()
$ws2.Select()
